$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.717.67"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "3.006.36"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'510.93"
$ws.Range("E5").Value = "  +6.63%  "
$ws.Range("D6").Value = "'138.94"
$ws.Range("E6").Value = "  +7.65%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +5.21%  "
$ws.Range("D9").Value = "'7.55"
$ws.Range("E9").Value = "  +9.77%  "
$ws.Range("E10").Value = "  +9.47%  "
$ws.Range("E11").Value = "  +4.09%  "
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").Value = "3.524.68"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").Value = "'25.78"
$ws.Range("E14").Value = "  +7.43%  "
$ws.Range("D15").Value = "'0.0000157"
$ws.Range("E15").Value = "  +14.11%  "
$ws.Range("D16").Value = "56.807.13"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "3.009.02"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "'5.97"
$ws.Range("E18").Value = "  +8.20%  "
$ws.Range("D19").Value = "'12.53"
$ws.Range("E19").Value = "  +6.38%  "
$ws.Range("D20").Value = "'7.86"
$ws.Range("E20").Value = "  +7.48%  "
$ws.Range("D21").Value = "'331.13"
$ws.Range("E21").Value = "  +7.72%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'0.483"
$ws.Range("E23").Value = "  +6.58%  "
$ws.Range("D24").Value = "'63.10"
$ws.Range("E24").Value = "  +6.18%  "
$ws.Range("E25").Value = "  +10.79%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "0.0₃0907"
$ws.Range("E27").Value = "  +9.09%  "
$ws.Range("D28").Value = "'6.71"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("D29").Value = "'7.09"
$ws.Range("E29").Value = "  +11.20%  "
$ws.Range("D30").Value = "'1.26"
$ws.Range("E30").Value = "  +8.96%  "
$ws.Range("E31").Value = "  +8.46%  "
$ws.Range("D32").Value = "'20.69"
$ws.Range("E32").Value = "  +8.74%  "
$ws.Range("D33").Value = "'154.40"
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("D34").Value = "'4.58"
$ws.Range("E34").Value = "  +7.09%  "
$ws.Range("D35").Value = "'5.68"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "'1.27"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").Value = "'0.0679"
$ws.Range("E37").Value = "  +6.82%  "
$ws.Range("D38").Value = "'24.13"
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("D39").Value = "3.043.07"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").Value = "'36.97"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").Value = "2.273.07"
$ws.Range("E43").Value = "  +8.21%  "
$ws.Range("E44").Value = "  +5.95%  "
$ws.Range("E45").Value = "  +4.95%  "
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "'1.98"
$ws.Range("E47").Value = "  +20.59%  "
$ws.Range("E48").Value = "  +6.96%  "
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("D50").Value = "'19.64"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("D51").Value = "'0.0873"
$ws.Range("E51").Value = "  +7.88%  "
